# Insert a new weekly price record at row 44, shifting all existing
# records from row 44 onward down by one row (44->45, ..., 148->149).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 44:148 down to make room for the new record.
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with the new record's data.
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 45133
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = 100112026
$ws.Range("G44").Value = "Haba"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 60
$ws.Range("K44").Value = 19000
$ws.Range("L44").Value = 19000
$ws.Range("M44").Value = 19000
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 760
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
